$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the title/name value in A17 (was "xyz", now "surya")
$ws.Range("A17").Value = "surya"

# Move the active selection from A18:XFD18 to A17 (single cell)
$ws.Range("A17").Select()
